$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B->C, old C->D)
$ws.Columns("B").Insert()

# New header for inserted column B1 -- copy formatting from the
# neighboring header cell (C1) so it matches the other header cells' style
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B1").Value = "segments"

# Move the segment names (currently in column A, rows 2-20) into new column B
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value2
}

# Set column A to numeric row index (0-based), rows 2-20
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
